$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value2 = 1999.7142
$ws.Cells.Item(70, 9).Value2 = 1999.6
$ws.Cells.Item(70, 11).Value2 = 5998.799999999999
$ws.Cells.Item(70, 13).Value2 = -5728.799999999999
$ws.Cells.Item(73, 8).Value2 = 1999.7142
$ws.Cells.Item(73, 9).Value2 = 1999.6
$ws.Cells.Item(73, 11).Value2 = 5998.799999999999
$ws.Cells.Item(73, 13).Value2 = -5062.799999999999
$ws.Cells.Item(76, 8).Value2 = 9999
$ws.Cells.Item(76, 10).Value2 = 9999
$ws.Cells.Item(76, 12).Value2 = 9999
$ws.Cells.Item(76, 14).Value2 = -10629
$ws.Cells.Item(79, 8).Value2 = 9999
$ws.Cells.Item(79, 10).Value2 = 9999
$ws.Cells.Item(79, 12).Value2 = 9999
$ws.Cells.Item(79, 14).Value2 = -12183
$ws.Cells.Item(86, 8).Value2 = 500001000
$ws.Cells.Item(86, 9).Value2 = 1000000000
$ws.Cells.Item(86, 11).Value2 = 1000000000
$ws.Cells.Item(86, 13).Value2 = -999998877
$ws.Cells.Item(89, 8).Value2 = 500001000
$ws.Cells.Item(89, 9).Value2 = 1000000000
$ws.Cells.Item(89, 11).Value2 = 5000000000
$ws.Cells.Item(89, 13).Value2 = -4999994384
$ws.Cells.Item(113, 8).Value2 = 5000
$ws.Cells.Item(113, 9).Value2 = 5000
$ws.Cells.Item(113, 11).Value2 = 5000
$ws.Cells.Item(113, 13).Value2 = -1746
$ws.Cells.Item(132, 8).Value2 = 6280.972
$ws.Cells.Item(132, 9).Value2 = 5153.4414
$ws.Cells.Item(132, 11).Value2 = 15460.3242
$ws.Cells.Item(132, 13).Value2 = -12930.3242
$ws.Cells.Item(133, 8).Value2 = 105980
$ws.Cells.Item(133, 10).Value2 = 105980
$ws.Cells.Item(133, 12).Value2 = 105980
$ws.Cells.Item(133, 14).Value2 = -116100
$ws.Cells.Item(137, 8).Value2 = 2506061.8
$ws.Cells.Item(137, 9).Value2 = 4167805.2
$ws.Cells.Item(137, 11).Value2 = 12503415.6
$ws.Cells.Item(137, 13).Value2 = -12500865.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value2 = 1449.8334
$ws.Cells.Item(2, 9).Value2 = 1166.6666
$ws.Cells.Item(2, 10).Value2 = 1733
$ws.Cells.Item(2, 11).Value2 = 1166.6666
$ws.Cells.Item(2, 12).Value2 = 1733
$ws.Cells.Item(2, 13).Value2 = -1053.6666
$ws.Cells.Item(2, 14).Value2 = -1959
$ws.Cells.Item(45, 8).Value2 = 31635.2
$ws.Cells.Item(45, 9).Value2 = 40819.184
$ws.Cells.Item(45, 11).Value2 = 40819.184
$ws.Cells.Item(45, 13).Value2 = -40442.184
$ws.Cells.Item(61, 8).Value2 = 3617.8948
$ws.Cells.Item(61, 9).Value2 = 2545.8438
$ws.Cells.Item(61, 11).Value2 = 2545.8438
$ws.Cells.Item(61, 13).Value2 = -2333.8438
$ws.Cells.Item(74, 8).Value2 = 224430.48
$ws.Cells.Item(74, 9).Value2 = 348448.44
$ws.Cells.Item(74, 10).Value2 = 3954.111
$ws.Cells.Item(74, 11).Value2 = 348448.44
$ws.Cells.Item(74, 12).Value2 = 3954.111
$ws.Cells.Item(74, 13).Value2 = -347574.44
$ws.Cells.Item(74, 14).Value2 = -5702.111
$ws.Cells.Item(77, 8).Value2 = 224430.48
$ws.Cells.Item(77, 9).Value2 = 348448.44
$ws.Cells.Item(77, 10).Value2 = 3954.111
$ws.Cells.Item(77, 11).Value2 = 1742242.2
$ws.Cells.Item(77, 12).Value2 = 19770.555
$ws.Cells.Item(77, 13).Value2 = -1737874.2
$ws.Cells.Item(77, 14).Value2 = -28506.555
$ws.Cells.Item(110, 8).Value2 = 3114.2222
$ws.Cells.Item(110, 9).Value2 = 1677.9445
$ws.Cells.Item(110, 10).Value2 = 5986.778
$ws.Cells.Item(110, 11).Value2 = 1677.9445
$ws.Cells.Item(110, 12).Value2 = 5986.778
$ws.Cells.Item(110, 13).Value2 = 367.0554999999999
$ws.Cells.Item(110, 14).Value2 = -10076.778
$ws.Cells.Item(116, 8).Value2 = 1449.8334
$ws.Cells.Item(116, 9).Value2 = 1166.6666
$ws.Cells.Item(116, 10).Value2 = 1733
$ws.Cells.Item(116, 11).Value2 = 1166.6666
$ws.Cells.Item(116, 12).Value2 = 1733
$ws.Cells.Item(116, 13).Value2 = 1127.3334
$ws.Cells.Item(116, 14).Value2 = -6321
$ws.Cells.Item(122, 8).Value2 = 3691.92
$ws.Cells.Item(122, 9).Value2 = 3491.1738
$ws.Cells.Item(122, 11).Value2 = 10473.5214
$ws.Cells.Item(122, 13).Value2 = -8023.5214
$ws.Cells.Item(132, 8).Value2 = 2318.7346
$ws.Cells.Item(132, 9).Value2 = 1369.7142
$ws.Cells.Item(132, 10).Value2 = 4691.2856
$ws.Cells.Item(132, 11).Value2 = 4109.142599999999
$ws.Cells.Item(132, 12).Value2 = 14073.8568
$ws.Cells.Item(132, 13).Value2 = -1579.142599999999
$ws.Cells.Item(132, 14).Value2 = -19133.8568
$ws.Cells.Item(136, 8).Value2 = 3617.8948
$ws.Cells.Item(136, 9).Value2 = 2545.8438
$ws.Cells.Item(136, 11).Value2 = 7637.5314
$ws.Cells.Item(136, 13).Value2 = -5087.5314
$ws.Cells.Item(139, 8).Value2 = 69999
$ws.Cells.Item(139, 10).Value2 = 69999
$ws.Cells.Item(139, 12).Value2 = 69999
$ws.Cells.Item(139, 14).Value2 = -80279

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value2 = 1449.8334
$ws.Cells.Item(3, 9).Value2 = 1166.6666
$ws.Cells.Item(3, 10).Value2 = 1733
$ws.Cells.Item(3, 11).Value2 = 1166.6666
$ws.Cells.Item(3, 12).Value2 = 1733
$ws.Cells.Item(3, 13).Value2 = -1052.6666
$ws.Cells.Item(3, 14).Value2 = -1961
$ws.Cells.Item(86, 8).Value2 = 3205.077
$ws.Cells.Item(86, 9).Value2 = 2997.3635
$ws.Cells.Item(86, 11).Value2 = 2997.3635
$ws.Cells.Item(86, 13).Value2 = -1874.3635
$ws.Cells.Item(89, 8).Value2 = 3205.077
$ws.Cells.Item(89, 9).Value2 = 2997.3635
$ws.Cells.Item(89, 11).Value2 = 14986.8175
$ws.Cells.Item(89, 13).Value2 = -9370.817499999999
$ws.Cells.Item(132, 8).Value2 = 83999.75
$ws.Cells.Item(132, 10).Value2 = 83999.75
$ws.Cells.Item(132, 12).Value2 = 83999.75
$ws.Cells.Item(132, 14).Value2 = -94119.75
$ws.Cells.Item(134, 8).Value2 = 5209.723
$ws.Cells.Item(134, 9).Value2 = 4853.4717
$ws.Cells.Item(134, 11).Value2 = 14560.4151
$ws.Cells.Item(134, 13).Value2 = -12025.4151

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 4839.3887
$ws.Cells.Item(31, 9).Value2 = 2947.55
$ws.Cells.Item(31, 10).Value2 = 7204.1875
$ws.Cells.Item(31, 11).Value2 = 2947.55
$ws.Cells.Item(31, 12).Value2 = 7204.1875
$ws.Cells.Item(31, 13).Value2 = -2652.55
$ws.Cells.Item(31, 14).Value2 = -7794.1875
$ws.Cells.Item(34, 8).Value2 = 4839.3887
$ws.Cells.Item(34, 9).Value2 = 2947.55
$ws.Cells.Item(34, 10).Value2 = 7204.1875
$ws.Cells.Item(34, 11).Value2 = 2947.55
$ws.Cells.Item(34, 12).Value2 = 7204.1875
$ws.Cells.Item(34, 13).Value2 = -2745.55
$ws.Cells.Item(34, 14).Value2 = -7608.1875
$ws.Cells.Item(41, 8).Value2 = 17702.5
$ws.Cells.Item(41, 9).Value2 = 5000
$ws.Cells.Item(41, 10).Value2 = 38873.332
$ws.Cells.Item(41, 11).Value2 = 5000
$ws.Cells.Item(41, 12).Value2 = 38873.332
$ws.Cells.Item(41, 13).Value2 = -4572
$ws.Cells.Item(41, 14).Value2 = -39729.332
$ws.Cells.Item(99, 8).Value2 = 11319.3
$ws.Cells.Item(99, 9).Value2 = 11354.777
$ws.Cells.Item(99, 10).Value2 = 11000
$ws.Cells.Item(99, 11).Value2 = 11354.777
$ws.Cells.Item(99, 12).Value2 = 11000
$ws.Cells.Item(99, 13).Value2 = -9856.777
$ws.Cells.Item(99, 14).Value2 = -13996
$ws.Cells.Item(126, 8).Value2 = 11319.3
$ws.Cells.Item(126, 9).Value2 = 11354.777
$ws.Cells.Item(126, 10).Value2 = 11000
$ws.Cells.Item(126, 11).Value2 = 34064.331
$ws.Cells.Item(126, 12).Value2 = 33000
$ws.Cells.Item(126, 13).Value2 = -31594.331
$ws.Cells.Item(126, 14).Value2 = -37940
$ws.Cells.Item(134, 8).Value2 = 2223.2341
$ws.Cells.Item(134, 9).Value2 = 1993.1951
$ws.Cells.Item(134, 10).Value2 = 3795.1667
$ws.Cells.Item(134, 11).Value2 = 5979.5853
$ws.Cells.Item(134, 12).Value2 = 11385.5001
$ws.Cells.Item(134, 13).Value2 = -3444.5853
$ws.Cells.Item(134, 14).Value2 = -16455.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value2 = 6224.875
$ws.Cells.Item(55, 9).Value2 = 1000
$ws.Cells.Item(55, 10).Value2 = 6971.2856
$ws.Cells.Item(55, 11).Value2 = 3000
$ws.Cells.Item(55, 12).Value2 = 20913.8568
$ws.Cells.Item(55, 13).Value2 = -2823
$ws.Cells.Item(55, 14).Value2 = -21267.8568
$ws.Cells.Item(121, 8).Value2 = 600319.8
$ws.Cells.Item(121, 9).Value2 = 666866.3
$ws.Cells.Item(121, 10).Value2 = 500500
$ws.Cells.Item(121, 11).Value2 = 2000598.9
$ws.Cells.Item(121, 12).Value2 = 1501500
$ws.Cells.Item(121, 13).Value2 = -1999288.9
$ws.Cells.Item(121, 14).Value2 = -1504120
$ws.Cells.Item(131, 8).Value2 = 19852.125
$ws.Cells.Item(131, 9).Value2 = 37382.5
$ws.Cells.Item(131, 11).Value2 = 112147.5
$ws.Cells.Item(131, 13).Value2 = -107107.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value2 = 100026.86
$ws.Cells.Item(70, 9).Value2 = 129412.875
$ws.Cells.Item(70, 11).Value2 = 129412.875
$ws.Cells.Item(70, 13).Value2 = -129142.875
$ws.Cells.Item(73, 8).Value2 = 100026.86
$ws.Cells.Item(73, 9).Value2 = 129412.875
$ws.Cells.Item(73, 11).Value2 = 129412.875
$ws.Cells.Item(73, 13).Value2 = -128476.875
$ws.Cells.Item(132, 8).Value2 = 7250.5
$ws.Cells.Item(132, 9).Value2 = 2431
$ws.Cells.Item(132, 11).Value2 = 7293
$ws.Cells.Item(132, 13).Value2 = -4763
$ws.Cells.Item(136, 8).Value2 = 13037.2
$ws.Cells.Item(136, 10).Value2 = 13037.2
$ws.Cells.Item(136, 12).Value2 = 39111.60000000001
$ws.Cells.Item(136, 14).Value2 = -44211.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 1355.75
$ws.Cells.Item(46, 9).Value2 = 499
$ws.Cells.Item(46, 10).Value2 = 1869.8
$ws.Cells.Item(46, 11).Value2 = 499
$ws.Cells.Item(46, 12).Value2 = 1869.8
$ws.Cells.Item(46, 13).Value2 = -311
$ws.Cells.Item(46, 14).Value2 = -2245.8
$ws.Cells.Item(122, 8).Value2 = 4781.8335
$ws.Cells.Item(122, 9).Value2 = 5258.6
$ws.Cells.Item(122, 11).Value2 = 15775.8
$ws.Cells.Item(122, 13).Value2 = -13325.8
$ws.Cells.Item(136, 8).Value2 = 4804.3125
$ws.Cells.Item(136, 9).Value2 = 4807.4546
$ws.Cells.Item(136, 11).Value2 = 14422.3638
$ws.Cells.Item(136, 13).Value2 = -11872.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value2 = 6566.6665
$ws.Cells.Item(62, 9).Value2 = 3850
$ws.Cells.Item(62, 11).Value2 = 3850
$ws.Cells.Item(62, 13).Value2 = -3226
$ws.Cells.Item(65, 8).Value2 = 6566.6665
$ws.Cells.Item(65, 9).Value2 = 3850
$ws.Cells.Item(65, 11).Value2 = 19250
$ws.Cells.Item(65, 13).Value2 = -16130
$ws.Cells.Item(136, 8).Value2 = 200011260
$ws.Cells.Item(136, 9).Value2 = 250001570
$ws.Cells.Item(136, 11).Value2 = 750004710
$ws.Cells.Item(136, 13).Value2 = -750002160
